$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PM25_compare")

# Row 4 (United States / usa) historical series N4:AL4 updated to 2023 EQUATES results.
$ws.Range("N4").Value = 177.66406130332712
$ws.Range("O4").Value = 180.70105380423871
$ws.Range("P4").Value = 181.4603019294666
$ws.Range("Q4").Value = 218.66346006563339
$ws.Range("R4").Value = 205.75624193675921
$ws.Range("S4").Value = 187.53428693128973
$ws.Range("T4").Value = 324.0159706741307
$ws.Range("U4").Value = 333.39420351694565
$ws.Range("V4").Value = 348.71810842842024
$ws.Range("W4").Value = 333.91580697897717
$ws.Range("X4").Value = 253.19102780849707
$ws.Range("Y4").Value = 253.50231953984053
$ws.Range("Z4").Value = 204.06679832
$ws.Range("AA4").Value = 202.26893871999999
$ws.Range("AB4").Value = 202.35532724000001
$ws.Range("AC4").Value = 202.35532724000001
$ws.Range("AD4").Value = 202.35657230999999
$ws.Range("AE4").Value = 202.28164340000001
$ws.Range("AF4").Value = 202.28164340000001
$ws.Range("AG4").Value = 202.28191641000001
$ws.Range("AH4").Value = 202.27903229
$ws.Range("AI4").Value = 202.76924987000001
$ws.Range("AJ4").Value = 202.82113631999999
$ws.Range("AK4").Value = 202.88598726000001
$ws.Range("AL4").Value = 202.85367769000001
